$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): add new column headers for season record (Wins, Losses, Ties)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the formatting from the existing header style (AC1) onto the new headers
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows 2-36: season record values (Wins, Losses, Ties)
for ($r = 2; $r -le 36; $r++) {
    $ws.Cells.Item($r, 30).Value = 101
    $ws.Cells.Item($r, 31).Value = 61
    $ws.Cells.Item($r, 32).Value = 0
}
